$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.533.99"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "2.536.63"
$ws.Range("E3").Value = "  -5.14%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.07"
$ws.Range("E5").Value = "  -3.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.80"
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.509"
$ws.Range("E8").Value = "  -2.57%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.166"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("B10").Value = "LidoStakedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D10").Value = "2.533.70"
$ws.Range("E10").Value = "  -5.12%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.344"
$ws.Range("E12").Value = "  -3.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.80"
$ws.Range("E13").Value = "  -3.36%  "
$ws.Range("D14").Value = "3.007.26"
$ws.Range("E14").Value = "  -5.15%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "70.378.02"
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000180"
$ws.Range("E16").Value = "  -2.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.03"
$ws.Range("E17").Value = "  -4.18%  "
$ws.Range("D18").Value = "2.535.99"
$ws.Range("E18").Value = "  -5.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.61"
$ws.Range("E19").Value = "  -4.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "360.29"
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.38"
$ws.Range("E21").Value = "  -9.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.94"
$ws.Range("E22").Value = "  -5.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.99"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.79"
$ws.Range("E25").Value = "  -3.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.10"
$ws.Range("E26").Value = "  -5.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.27"
$ws.Range("E27").Value = "  -4.75%  "
$ws.Range("D28").Value = "2.677.51"
$ws.Range("E28").Value = "  -4.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").Value = "0.0₃0925"
$ws.Range("E30").Value = "  -4.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.89"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "484.09"
$ws.Range("E32").Value = "  -3.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.28"
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.76"
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +6.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.29"
$ws.Range("E37").Value = "  -3.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.63"
$ws.Range("E38").Value = "  -4.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.82"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.31"
$ws.Range("E41").Value = "  -4.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.67"
$ws.Range("E42").Value = "  -5.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.74"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.45"
$ws.Range("E44").Value = "  -4.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.319"
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.41"
$ws.Range("E46").Value = "  -2.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "144.44"
$ws.Range("E47").Value = "  -7.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.53"
$ws.Range("E48").Value = "  -4.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.527"
$ws.Range("E49").Value = "  -5.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.63"
$ws.Range("E50").Value = "  -6.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.594"
$ws.Range("E51").Value = "  -1.81%  "
